$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '70.304.26'
$ws.Range('E2').Value2 = '  +4.62%  '
$ws.Range('D3').Value2 = '3.629.64'
$ws.Range('E3').Value2 = '  +4.40%  '
$ws.Range('E4').Value2 = '  +0.01%  '
$ws.Range('D5').Value2 = '''591.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +1.29%  '
$ws.Range('D6').Value2 = '''193.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  +4.24%  '
$ws.Range('D7').Value2 = '''0.644'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +1.78%  '
$ws.Range('D8').Value2 = '3.623.64'
$ws.Range('E8').Value2 = '  +4.41%  '
$ws.Range('D10').Value2 = '''0.181'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  +4.14%  '
$ws.Range('D11').Value2 = '''0.669'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  +2.95%  '
$ws.Range('D12').Value2 = '''58.39'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  +3.79%  '
$ws.Range('D13').Value2 = '''0.0000291'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  +4.17%  '
$ws.Range('D14').Value2 = '''9.94'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +5.44%  '
$ws.Range('D15').Value2 = '4.215.05'
$ws.Range('E15').Value2 = '  +4.31%  '
$ws.Range('D16').Value2 = '''19.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  +5.41%  '
$ws.Range('D17').Value2 = '3.630.19'
$ws.Range('E17').Value2 = '  +4.07%  '
$ws.Range('D18').Value2 = '70.306.19'
$ws.Range('E18').Value2 = '  +4.50%  '
$ws.Range('D19').Value2 = '''12.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +4.79%  '
$ws.Range('E20').Value2 = '  +1.63%  '
$ws.Range('E21').Value2 = '  +4.48%  '
$ws.Range('D22').Value2 = '''489.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +0.50%  '
$ws.Range('D23').Value2 = '''19.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  +14.77%  '
$ws.Range('E24').Value2 = '  -0.82%  '
$ws.Range('E25').Value2 = '  +0.37%  '
$ws.Range('D26').Value2 = '''91.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value2 = '  +1.72%  '
$ws.Range('D27').Value2 = '''3.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  +7.77%  '
$ws.Range('D28').Value2 = '''11.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +5.15%  '
$ws.Range('D29').Value2 = '''9.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  +5.17%  '
$ws.Range('D30').Value2 = '''32.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +4.95%  '
$ws.Range('D31').Value2 = '''7.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  +9.83%  '
$ws.Range('E32').Value2 = '  +8.46%  '
$ws.Range('D33').Value2 = '''627.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  +3.66%  '
$ws.Range('E34').Value2 = '  +4.30%  '
$ws.Range('D35').Value2 = '''65.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  +2.65%  '
$ws.Range('D36').Value2 = '''40.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  +11.47%  '
$ws.Range('E37').Value2 = '  +7.03%  '
$ws.Range('D38').Value2 = '0.0₃0822'
$ws.Range('E38').Value2 = '  +8.37%  '
$ws.Range('B39').Value2 = 'Kaspa'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value2 = '''0.147'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  -1.47%  '
$ws.Range('B40').Value2 = 'Dai'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value2 = '''1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  +0.02%  '
$ws.Range('D41').Value2 = '''3.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +0.65%  '
$ws.Range('D42').Value2 = '3.295.62'
$ws.Range('E42').Value2 = '  +0.88%  '
$ws.Range('D43').Value2 = '''3.19'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  +9.39%  '
$ws.Range('D44').Value2 = '''2.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  +12.69%  '
$ws.Range('D45').Value2 = '''0.0454'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +5.63%  '
$ws.Range('D46').Value2 = '''2.88'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +2.95%  '
$ws.Range('E47').Value2 = '  +0.92%  '
$ws.Range('E48').Value2 = '  +2.49%  '
$ws.Range('D49').Value2 = '''9.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value2 = '  +4.84%  '
$ws.Range('D50').Value2 = '''3.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +1.96%  '
$ws.Range('E51').Value2 = '  -0.12%  '
